$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "d" to G8 (Visualisasi Data) and "td" to G9, G10 (Twittr API, Rvest Scraping)
$ws.Range("G8").Value = "d"
$ws.Range("G9").Value = "td"
$ws.Range("G10").Value = "td"

# Narrow column G width to fit the new short labels (stored width 6.5)
$ws.Range("G:G").ColumnWidth = 5.666666666666667

# Move the active cell selection from I15 to H15
$ws.Range("H15").Select()
